# The presentation originally ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Integral" (used by the Slide Master)
#   ppt/theme/theme2.xml -> clrScheme "Office"   (used only by the Notes Master)
#
# The target edit swaps the two themes' colour schemes, so that the theme
# reachable from the Slide Master (theme1.xml) ends up carrying the
# "Office" colour values, while the theme used by the Notes Master
# (theme2.xml) ends up carrying the "Integral" colour values.
#
# PowerPoint's ThemeColorScheme object model addresses theme colours as a
# flat, 12-slot, 1-based collection in a fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# and each ThemeColor.RGB is a VBA-style 0xBBGGRR colour long (the classic
# COLORREF byte order), not 0xRRGGBB.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target values are the "Office" theme's RGB (0xRRGGBB) colours, in
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.
$officeRgbHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeRgbHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)

    # VBA/COM RGB long is 0x00BBGGRR.
    $bgrLong = ($b * 65536) + ($g * 256) + $r

    $colorScheme.Colors($i).RGB = $bgrLong
}

Write-Output "Swapped theme colour scheme onto the Slide Master's theme (Office colours applied)."
